$d = $word.ActiveDocument

# --- Edit 1: "Помощь Общества будет направлена для приобретения жилого
# помещения квартира/дом (нужное подчеркнуть) у лиц, ..." -->
# "... жилого помещения (${JP_TYPE}) у лиц, ..."

# 1a. Drop the trailing "квартира/дом" from the first run, leaving the
#     sentence ending in a single trailing space.
$r1 = $d.Content
$r1.Find.Execute("квартира/дом", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Text = ""

# 1b. The old " (нужное подчеркнуть" run (underlined) becomes "($" with no
#     underline.
$r2 = $d.Content
$r2.Find.Execute(" (нужное подчеркнуть", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Text = "(`$"
$r2.Font.Underline = 0

# 1c. Insert the bold placeholder "{JP_TYPE}" right after it (still
#     un-underlined), ahead of the untouched ") у лиц, ..." run.
$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter("{JP_TYPE}")
$r3.Font.Bold = 1
$r3.Font.Underline = 0

# --- Edit 2: "... я и члены моей семьи не участвовали." -->
# "... я и члены моей семьи ${IS_PARTICIPATE}."

$r4 = $d.Content
$r4.Find.Execute("не участвовали", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4.Text = "`${IS_PARTICIPATE}"
$r4.Font.Bold = 1
